# Update parameter inputs for new scenarios
# (mirrors the "update parameter inputs for new scenarios" commit)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# ---------------------------------------------------------------------------
# Row 6 - e / SelfTestUnit : Georgia/Kenya/China/Vietnam cost bumped 5 -> 5.63
# ---------------------------------------------------------------------------
$ws.Range("F6").Value = 5.63
$ws.Range("K6").Value = 5.63
$ws.Range("P6").Value = 5.63
$ws.Range("U6").Value = 5.63

# ---------------------------------------------------------------------------
# Row 7 - f / StandardAntibody : replace formulas with literal rounded values
# ---------------------------------------------------------------------------
$ws.Range("F7").Value = 32.76
$ws.Range("K7").Value = 42.87
$ws.Range("P7").Value = 5.07

# ---------------------------------------------------------------------------
# Row 15 - n / NewSelfTest
# ---------------------------------------------------------------------------
$ws.Range("F15").Value = 0.09
$ws.Range("K15").Value = 0.36
$ws.Range("P15").Value = 0.07
$ws.Range("U15").Value = 0.36

# ---------------------------------------------------------------------------
# Row 16 - o / SelfTestInstead
# ---------------------------------------------------------------------------
$ws.Range("F16").Value = 0.1
$ws.Range("K16").Value = 0.1
$ws.Range("P16").Value = 0.1
$ws.Range("U16").Value = 0.1

# ---------------------------------------------------------------------------
# Row 24 - w / SelfTestSuccess
# ---------------------------------------------------------------------------
$ws.Range("F24").Value = 0.97
$ws.Range("K24").Value = 0.97
$ws.Range("P24").Value = 0.97
$ws.Range("U24").Value = 0.97

# ---------------------------------------------------------------------------
# Row 28 - z1 / PercentReportResultsPos
# ---------------------------------------------------------------------------
$ws.Range("K28").Value = 0.65
$ws.Range("U28").Value = 0.65

# ---------------------------------------------------------------------------
# Row 30 - z3 / PercentReportResultsUnknown
# ---------------------------------------------------------------------------
$ws.Range("K30").Value = 0.65
$ws.Range("U30").Value = 0.65

# ---------------------------------------------------------------------------
# Rows 32-36 : rename the OraQuick/Retest parameters to Self-test /
# Facility-based naming, introduce a new "Facility based" intervention row
# ---------------------------------------------------------------------------

# Row 32 - sens : "OraQuick sensitivity" -> "Self-test sensitivity"; drop the
# now-unused ParameterType ("Intervention") cell entirely
$ws.Range("B32").Value = "Self-test sensitivity"
$ws.Range("C32").ClearContents()

# Row 33 - spec : "OraQuick specificity" -> "Self-test specificity"
$ws.Range("B33").Value = "Self-test specificity"

# Row 35 - sens2 : "Retest sensitivity" -> "Facility-based sensitivity",
# tag it with the new "Facility based" ParameterType, and drop sensitivity
# from 100% to 95%
$ws.Range("B35").Value = "Facility-based sensitivity"
$ws.Range("C35").Value = "Facility based"
$ws.Range("F35").Value = 0.95
$ws.Range("K35").Value = 0.95
$ws.Range("P35").Value = 0.95
$ws.Range("U35").Value = 0.95

# Row 36 - spec2 : "Retest specificity" -> "Facility-based specificity"
$ws.Range("B36").Value = "Facility-based specificity"

# ---------------------------------------------------------------------------
# Make "Parameters" the active/selected sheet with a frozen header
# row + label column, scrolled to the Vietnam columns, mirroring the
# workbook's new default view (previously "Self-report" was the active tab)
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("C2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("U36").Select()
